# The document starts as a single paragraph:
#   "This is a test Microsoft Word document for github."
# wrapped by the hidden "_GoBack" bookmark (spanning the whole run text),
# with a proofErr spell-check marker around "github".
#
# Target:
#   Run 1: "This is a test Microsoft Word document for github"
#   Run 2: " and bitbucket"
#   (bookmark "_GoBack" now collapsed, sitting right before the final period)
#   Run 3: "."
#
# Rebuild the paragraph's runs from scratch (typing fresh text onto an
# emptied range keeps each InsertAfter() call as its own run instead of
# Word's usual "merge into the adjacent same-format run" behavior), then
# relocate the _GoBack bookmark to its new collapsed position.

$d = $word.ActiveDocument

$para = $d.Paragraphs(1)
$paraEnd = $para.Range.End

# Wipe the paragraph's text (leaves the bookmark collapsed at position 0,
# and drops the now-stale proofErr spell-check markers).
$body = $d.Range(0, $paraEnd - 1)
$body.Delete()

# Run 1
$r1 = $d.Range(0, 0)
$r1.InsertAfter("This is a test Microsoft Word document for github")

# Run 2
$p1End = $d.Content.End - 1
$r2 = $d.Range($p1End, $p1End)
$r2.InsertAfter(" and bitbucket")

# Run 3 (the trailing period)
$p2End = $d.Content.End - 1
$r3 = $d.Range($p2End, $p2End)
$r3.InsertAfter(".")

# Move the _GoBack bookmark so it collapses right before the final period.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()
$periodPos = $d.Content.End - 2
$bmRange = $d.Range($periodPos, $periodPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
